# Daily attendance processing - 2026-01-17 14:04:15
#
# The "Recorded By" column (G) stores a comma-separated list of users who
# recorded/touched a session. For the rows below, "System" was the first
# entry in the list; this pass normalizes those rows by moving the
# "System" marker from the front of the list to the back.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows in column G whose "Recorded By" list needs "System" moved from the
# first position to the last position.
$rowsToFix = @(2,3,4,5,6,7,8,28,29,30,31,32,33,34,54,55,56,57,58,59,60,80,81,82,106,107,108,132,133,134)

foreach ($row in $rowsToFix) {
    $cell = $ws.Cells.Item($row, 7)  # Column G = 7
    $current = [string]$cell.Value2

    $parts = $current -split ',\s*'
    if ($parts.Length -gt 1) {
        $rotated = ($parts[1..($parts.Length - 1)] + $parts[0]) -join ', '
        $cell.Value2 = $rotated
    }
}
